$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: restyle to the "bordered" row look (same style pattern as
#     rows 3/6/8: A/B = style 8, C/D/E = style 9) and clear B10:E10 so only
#     the filename in A10 remains. Copy formats-only from row 8, which
#     already carries that exact style pattern, so no new style records
#     are needed.
$ws.Range("A8:E8").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# --- Row 11: new data row, same style pattern as row 7 (style 4 / 5).
$ws.Range("A7:E7").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# --- Row 12: continuation row (no filename cell), same style pattern as
#     row 5 (style 4 / 5), copied from columns B:E only.
$ws.Range("B5:E5").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)

# --- Populate the new cell text in the exact order the original author
#     typed it (this controls the order new entries land in
#     xl/sharedStrings.xml, which the target file depends on).
$ws.Range("C11").Value = " The planet will be paralyzed if\nsomething\'s not done?"
$ws.Range("C12").Value = " We\'d offer to help...[K]but being\nmotivated isn\'t really our thing…"
$ws.Range("A11").Value = "SCRIPT/G01P03A/um2204.ssb"
$ws.Range("D11").Value = " Если что-то не сделать, то\nпланету парализует?"
$ws.Range("D12").Value = " Мы бы с радостью помогли...[K]\nНо мотивация это не про нас..."
$ws.Range("E11").Value = " Åòìé œóï-óï îå òäåìàóû, óï\nðìàîåóô ðàñàìéèôåó?"
$ws.Range("E12").Value = " Íú áú ò ñàäïòóûý ðïíïãìé...[K]\nÎï íïóéâàøéÿ üóï îå ðñï îàò..."
$ws.Range("B11").Value = 67
$ws.Range("B12").Value = 70

# --- Row heights for the two new rows (43.2 / 31.8), matching the rest of
#     the sheet's wrapped-text rows.
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 31.8

# --- Final selection, as recorded in the saved sheet view.
$ws.Range("D10").Select()
